# DateRange update: insert the latest day (12/23/2020) at the top of the
# TSA flying-stats table, shifting all existing rows down by one and
# dropping the oldest row (6/6/2020) that falls off the bottom of the
# fixed 200-row window (A1:C200).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing 200 rows of data down by one row.
$ws.Rows("1:1").Insert()

# Drop the oldest row that was pushed past row 200.
$ws.Rows("201:201").Delete()

# The sheet stores every value (dates and numbers alike) as plain text,
# so force Text format on the new row before writing values - otherwise
# Excel would auto-convert "12/23/2020" into a date serial number and
# "1,191,123" into a numeric value.
$newRow = $ws.Range("A1:C1")
$newRow.NumberFormat = "@"

$ws.Range("A1").Value = "12/23/2020"
$ws.Range("B1").Value = "1,191,123"
$ws.Range("C1").Value = "1,937,235"
